$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2, pushing current rows 2-3 down to 3-4
$ws.Rows.Item(2).Insert()
# The inserted row picks up formatting from the row above (header); reset it
# to plain/unstyled so it matches the other data rows.
$ws.Rows.Item(2).ClearFormats()

# Copy the (now-shifted) original row 2 data (now in row 3) into new row 2
for ($col = 1; $col -le 20; $col++) {
    $ws.Cells.Item(2, $col).Value = $ws.Cells.Item(3, $col).Value()
}

# Column D (Fecha) carries the date number format style used by the other rows
$ws.Cells.Item(2, 4).NumberFormat = $ws.Cells.Item(3, 4).NumberFormat

# Apply the new/changed values for row 2 per the diff
$ws.Cells.Item(2, 4).Value = 45243   # D2 Fecha
$ws.Cells.Item(2, 13).Value = 50     # M2 Volumen
$ws.Cells.Item(2, 14).Value = 35000  # N2 Precio minimo
$ws.Cells.Item(2, 15).Value = 35000  # O2 Precio maximo
$ws.Cells.Item(2, 16).Value = 35000  # P2 Precio promedio ponderado
$ws.Cells.Item(2, 18).Value = "Región Metropolitana"  # R2 Origen
$ws.Cells.Item(2, 19).Value = 7000   # S2 Precio $/Kg
